$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores every Price/Volume cell as literal text (even values
# that look numeric, e.g. "1.00", "585.55"). Excel auto-converts such literals to
# numbers on assignment, so cells whose new text would otherwise be re-typed as a
# number are first switched to the Text number format to preserve the exact string.
$ws.Range("D2").Value = "70.924.42"
$ws.Range("E2").Value = "  +2.16%  "
$ws.Range("D3").Value = "3.572.21"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.55"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "190.63"
$ws.Range("E6").Value = "  +2.47%  "
$ws.Range("E7").Value = "  +1.85%  "
$ws.Range("D8").Value = "3.562.46"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +15.44%  "
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.72"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("E13").Value = "  +4.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.54"
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("D15").Value = "4.136.26"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "70.851.58"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "12.85"
$ws.Range("E17").Value = "  +4.28%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.576.64"
$ws.Range("E19").Value = "  +1.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "569.61"
$ws.Range("E20").Value = "  +3.89%  "
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.71"
$ws.Range("E23").Value = "  -4.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.60"
$ws.Range("E24").Value = "  +3.08%  "
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "94.39"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.25"
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.32"
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.68"
$ws.Range("E30").Value = "  +2.34%  "
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.35"
$ws.Range("E32").Value = "  -2.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.117"
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.80"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.80"
$ws.Range("E35").Value = "  +23.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.29"
$ws.Range("E36").Value = "  +6.10%  "
$ws.Range("E37").Value = "  +2.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "532.16"
$ws.Range("E38").Value = "  -3.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.47"
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "3.651.72"
$ws.Range("E40").Value = "  +10.12%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0802"
$ws.Range("E41").Value = "  +4.58%  "
$ws.Range("B42").Value = "Dai"
$ws.Range("C42").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.139"
$ws.Range("E43").Value = "  +4.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.46"
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0472"
$ws.Range("E45").Value = "  +5.77%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.95"
$ws.Range("E46").Value = "  -1.44%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.46"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  +4.10%  "
$ws.Range("E49").Value = "  +3.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("E51").Value = "  +7.55%  "

